# About us page automation
$wb = $excel.ActiveWorkbook

# Rename "Sheet2" to "data"
$dataSheet = $wb.Worksheets.Item("Sheet2")
$dataSheet.Name = "data"

# Populate the data sheet with the "About us" Lorem Ipsum content
$loremText = "There are many variations of passages of Lorem Ipsum available, but the majority have suffered alteration in some form, by injected humour, or randomised words which don't look even slightly believable. If you are going to use a passage of Lorem Ipsum, you need to be sure there isn't anything embarrassing hidden in the middle of text. All the Lorem Ipsum generators on the Internet tend to repeat predefined chunks as necessary, making this the first true generator on the Internet."

$dataSheet.Columns.Item(1).ColumnWidth = 50.1667
$dataSheet.Range("A1").Value = $loremText
$dataSheet.Range("A1").VerticalAlignment = -4108
$dataSheet.Range("A1").WrapText = $true
$dataSheet.Rows.Item(1).RowHeight = 150

# Make the "data" sheet the active / selected tab
$dataSheet.Activate()
